$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new task rows added at the bottom of the tracking table, matching the
# date format already used by the other "DATE" column entries.
$dateFormat = $ws.Range("C24").NumberFormat

$ws.Range("B25").Value = "preparation travaux a venir "
$ws.Range("B27").Value = "partie model (dossier data) les controlleurs avec requete sql 2"
$ws.Range("B26").Value = "partie model (dossier data) les controlleurs avec requete sql 1"

$ws.Range("C25").NumberFormat = $dateFormat
$ws.Range("C25").Value = 42338
$ws.Range("D25").Value = 1

$ws.Range("C26").NumberFormat = $dateFormat
$ws.Range("C26").Value = 42339
$ws.Range("D26").Value = 4

$ws.Range("C27").NumberFormat = $dateFormat
$ws.Range("C27").Value = 42341
$ws.Range("D27").Value = 3.5

# Update the visible scroll position / selection to match the author's final view
try {
    $ws.Application.ActiveWindow.ScrollRow = 10
    $ws.Application.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("D27").Select()
